$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (target stored width 15.42578125 chars)
$ws.Columns.Item(1).ColumnWidth = 14.592447916666666
$ws.Columns.Item(2).ColumnWidth = 14.592447916666666

# Update cell values
$ws.Range("A1").Value = 0.005034448180838506
$ws.Range("B1").Value = -0.0050344482006422601

$ws.Range("A2").Value = -0.0065336135539141877
$ws.Range("B2").Value = 0.0065336135412281671

$ws.Range("A3").Value = -0.034464233493570945
$ws.Range("B3").Value = 0.034464233480775743

$ws.Range("A4").Value = 0.046828604005686648
$ws.Range("B4").Value = -0.046828604016710733
